# Generate Report for Handoff
# The file "333085f4-f54c-4a78-8c04-096cc915e0fe" (row 3 on every sheet) is
# being handed off again: its status flips from "Handed back: in sync with
# en-US" to "Ready for handoff", and a new (later) handoff datetime is
# recorded for each locale.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# ----- Overview sheet -----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ----- zh-cn sheet -----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = $newStatus
$zhcn.Range("D3").Value = "2016-03-03 09:07:25"

# ----- de-de sheet -----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = $newStatus
$dede.Range("D3").Value = "2016-03-03 09:07:36"
